# Coursework.pptx edit: "Completed V0.1. Need to refine controls and minimap"
#
#  1. Slide 5 ("Bugs and Todo"): merge the two runs of the last bullet
#     ("Optimisation " + "through profiling") into a single run reading
#     "Optimisation through profiling".
#  2. Append a brand-new slide (Title + Content layout) after slide 5 with
#     a debugger-watch-style note about SearchAndDestoy.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 5: merge "Optimisation " + "through profiling" into one run.
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$bullets = $slide5.Shapes.Item(2).TextFrame.TextRange
$fullText = $bullets.Text
$needle = "Optimisation "
$start = $fullText.IndexOf($needle) + 1          # COM ranges are 1-based
$len = "Optimisation through profiling".Length
$bullets.Characters($start, $len).Text = "Optimisation through profiling"

# ---------------------------------------------------------------------
# 2. New slide, appended at the end, "Title and Content" layout (same
#    layout used by the other content slides in the deck).
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder is left blank (untouched) - just stamp the language
# like the rest of the deck uses.
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "x"
$titleRange.LanguageID = "en-GB"
$titleRange.Text = ""

# Content placeholder: three runs followed by a trailing blank paragraph.
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.LanguageID = "en-GB"
$body.Text = "-`t`tthis`t0x0bb9aa58 {...}`t"
$run2 = $body.InsertAfter("SearchAndDestoy")
$run3 = $run2.InsertAfter(" *")
$run4 = $run3.InsertAfter("`r")
